# Updated cryptos list on Thu Jun 13 13:44:14 UTC 2024 with GitHub Actions
#
# Note: several Price cells (column D) hold numeric-looking text such as
# "10.00" or "0.489" that must stay literal text (not be coerced to a
# number, which would drop the formatting). A leading apostrophe forces
# Excel to treat the assigned value as text, exactly like typing it by
# hand in the UI; the apostrophe itself is not stored in the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.243.63'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '3.519.63'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''610.14'
$ws.Range('E5').Value = '  -2.80%  '
$ws.Range('D6').Value = '''151.68'
$ws.Range('E6').Value = '  -4.84%  '
$ws.Range('D7').Value = '3.518.99'
$ws.Range('E7').Value = '  -2.76%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.489'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').Value = '''7.58'
$ws.Range('E11').Value = '  +4.87%  '
$ws.Range('D12').Value = '''0.432'
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('E13').Value = '  -3.33%  '
$ws.Range('D14').Value = '''32.23'
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('D15').Value = '4.114.07'
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').Value = '3.515.74'
$ws.Range('E16').Value = '  -2.52%  '
$ws.Range('D17').Value = '68.109.68'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '''6.54'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = '''15.52'
$ws.Range('E20').Value = '  -2.63%  '
$ws.Range('D21').Value = '''10.00'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').Value = '''451.49'
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('D24').Value = '''79.18'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('D25').Value = '3.656.63'
$ws.Range('E25').Value = '  -2.84%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -8.20%  '
$ws.Range('D28').Value = '''8.76'
$ws.Range('E28').Value = '  -5.17%  '
$ws.Range('D29').Value = '''10.02'
$ws.Range('E29').Value = '  -5.63%  '
$ws.Range('D30').Value = '''1.68'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('D31').Value = '''2.53'
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('E32').Value = '  -3.37%  '
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = '''25.73'
$ws.Range('E34').Value = '  -2.77%  '
$ws.Range('D35').Value = '''6.24'
$ws.Range('E35').Value = '  -5.08%  '
$ws.Range('D36').Value = '''1.86'
$ws.Range('E36').Value = '  -5.46%  '
$ws.Range('D37').Value = '3.512.43'
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('D38').Value = '''8.07'
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = '''2.31'
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''177.49'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '''0.998'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').Value = '''0.0910'
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('E44').Value = '  -2.68%  '
$ws.Range('D45').Value = '''31.38'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').Value = '''0.901'
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('D47').Value = '''47.04'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('D48').Value = '''1.32'
$ws.Range('E48').Value = '  -4.06%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '''2.54'
$ws.Range('E49').Value = '  -9.55%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '''7.66'
$ws.Range('E50').Value = '  -1.64%  '
